$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("summary")
$wsModelFit = $wb.Worksheets.Item("model_fit")

# Update the "summary" sheet (irt_dich table)
$wsSummary.Range("D2").Value = 1440
$wsSummary.Range("E2").Value = 82.22
$wsSummary.Range("F2").Value = -1.9
$wsSummary.Range("I2").Value = 0.26
$wsSummary.Range("J2").Value = 0.27
$wsSummary.Range("L2").Value = 1.04
$wsSummary.Range("D3").Value = 1431
$wsSummary.Range("E3").Value = 79.11
$wsSummary.Range("F3").Value = -1.66
$wsSummary.Range("H3").Value = 1.01
$wsSummary.Range("I3").Value = 0.39
$wsSummary.Range("J3").Value = 0.29
$wsSummary.Range("L3").Value = 1.06
$wsSummary.Range("D4").Value = 1428
$wsSummary.Range("E4").Value = 70.66
$wsSummary.Range("F4").Value = -1.11
$wsSummary.Range("I4").Value = -0.08
$wsSummary.Range("J4").Value = 0.35
$wsSummary.Range("L4").Value = 1.17
$wsSummary.Range("D5").Value = 1429
$wsSummary.Range("E5").Value = 67.32
$wsSummary.Range("F5").Value = -0.91
$wsSummary.Range("I5").Value = 1.07
$wsSummary.Range("D6").Value = 1425
$wsSummary.Range("E6").Value = 63.79
$wsSummary.Range("F6").Value = -0.71
$wsSummary.Range("I6").Value = -1.17
$wsSummary.Range("K6").Value = 0.03
$wsSummary.Range("L6").Value = 1.37
$wsSummary.Range("D7").Value = 1397
$wsSummary.Range("E7").Value = 53.26
$wsSummary.Range("H7").Value = 0.99
$wsSummary.Range("I7").Value = -0.45
$wsSummary.Range("D8").Value = 1388
$wsSummary.Range("E8").Value = 52.45
$wsSummary.Range("F8").Value = -0.12
$wsSummary.Range("H8").Value = 0.99
$wsSummary.Range("I8").Value = -0.34
$wsSummary.Range("J8").Value = 0.38
$wsSummary.Range("L8").Value = 1.2
$wsSummary.Range("D9").Value = 1349
$wsSummary.Range("E9").Value = 45.74
$wsSummary.Range("F9").Value = 0.23
$wsSummary.Range("H9").Value = 0.99
$wsSummary.Range("I9").Value = -0.5
$wsSummary.Range("J9").Value = 0.39
$wsSummary.Range("K9").Value = 0.03
$wsSummary.Range("L9").Value = 1.23
$wsSummary.Range("D10").Value = 1308
$wsSummary.Range("E10").Value = 39.45
$wsSummary.Range("F10").Value = 0.53
$wsSummary.Range("H10").Value = 1
$wsSummary.Range("I10").Value = -0.01
$wsSummary.Range("J10").Value = 0.37
$wsSummary.Range("K10").Value = 0.03
$wsSummary.Range("L10").Value = 1.15
$wsSummary.Range("D11").Value = 1246
$wsSummary.Range("E11").Value = 36.84
$wsSummary.Range("F11").Value = 0.69
$wsSummary.Range("G11").Value = 0.07
$wsSummary.Range("H11").Value = 0.98
$wsSummary.Range("I11").Value = -0.72
$wsSummary.Range("K11").Value = 0.03
$wsSummary.Range("L11").Value = 1.29
$wsSummary.Range("D12").Value = 1178
$wsSummary.Range("E12").Value = 34.04
$wsSummary.Range("F12").Value = 0.83
$wsSummary.Range("H12").Value = 1.02
$wsSummary.Range("I12").Value = 0.52
$wsSummary.Range("L12").Value = 1.09
$wsSummary.Range("D13").Value = 946
$wsSummary.Range("E13").Value = 22.3
$wsSummary.Range("F13").Value = 1.56
$wsSummary.Range("G13").Value = 0.09
$wsSummary.Range("H13").Value = 1.06
$wsSummary.Range("I13").Value = 1.2
$wsSummary.Range("J13").Value = 0.27
$wsSummary.Range("K13").Value = 0.04
$wsSummary.Range("L13").Value = 0.88
$wsSummary.Range("D14").Value = 460
$wsSummary.Range("E14").Value = 21.52
$wsSummary.Range("F14").Value = 1.63
$wsSummary.Range("H14").Value = 0.96
$wsSummary.Range("I14").Value = -0.52
$wsSummary.Range("J14").Value = 0.35
$wsSummary.Range("L14").Value = 1.41

# Update the "model_fit" sheet
$wsModelFit.Range("D2").Value = 18801
$wsModelFit.Range("E2").Value = 18829
$wsModelFit.Range("F2").Value = 18903
$wsModelFit.Range("G2").Value = 0.701
$wsModelFit.Range("H2").Value = 0.607
$wsModelFit.Range("D3").Value = 18784
$wsModelFit.Range("E3").Value = 18836
$wsModelFit.Range("F3").Value = 18974
$wsModelFit.Range("G3").Value = 0.703
$wsModelFit.Range("H3").Value = 0.604
